# Adds a new "event" column to test.xlsx (win/loss for competitive trials,
# rewarded/omission for reward trials), pushing the old "competition_closeness"
# column one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at M; existing M (competition_closeness) and
# everything right of it shifts one column right (M->N, N->O), and the
# data validation / dimension / formulas referencing those columns are
# updated automatically by Excel.
$ws.Columns("M").Insert()

# Header for the freshly inserted column.
$ws.Range("M1").Value = "event"

# Rows 2-21 describe head-to-head competition trials where column L holds a
# numeric condition code: 1.1 = subject 1 favored (a "win"), 1.4 = subject 2
# favored (a "loss").
for ($r = 2; $r -le 21; $r++) {
    $cond = $ws.Cells.Item($r, 12).Value()
    if ($cond -eq 1.1) {
        $ws.Cells.Item($r, 13).Value = "win"
    } else {
        $ws.Cells.Item($r, 13).Value = "loss"
    }
}

# Rows 22-39 are reward trials; column L already says "rewarded"/"omission" -
# mirror that same text into the new event column.
for ($r = 22; $r -le 39; $r++) {
    $cond = $ws.Cells.Item($r, 12).Value()
    $ws.Cells.Item($r, 13).Value = $cond
}

# Column widths that were hand-tuned for readability (values below are the
# ColumnWidth inputs that reproduce the saved pixel widths).
$ws.Columns("A").ColumnWidth = 10
$ws.Columns("B").ColumnWidth = 14.5
$ws.Columns("D").ColumnWidth = 62.333333333333336
$ws.Columns("E").ColumnWidth = 71.83333333333333
$ws.Columns("H").ColumnWidth = 32
$ws.Columns("K").ColumnWidth = 38

# Leave the selection where the editor last left it.
$ws.Range("O25").Select()
